$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 120 (shifts existing rows 120-177 down to 121-178)
$ws.Rows(120).Insert()

# Populate the new row 120 with data. Columns A,B,C,E,F,G,H,I,N,O,Q,R keep the
# same values as the row directly below (the former row 120, now row 121),
# while D,J,K,L,M,P get the new week's values.
$ws.Range("A120").Value = $ws.Range("A121").Value2
$ws.Range("B120").Value = $ws.Range("B121").Value2
$ws.Range("C120").Value = $ws.Range("C121").Value2
$ws.Range("D120").Value = 44455
$ws.Range("E120").Value = $ws.Range("E121").Value2
$ws.Range("F120").Value = $ws.Range("F121").Value2
$ws.Range("G120").Value = $ws.Range("G121").Value2
$ws.Range("H120").Value = $ws.Range("H121").Value2
$ws.Range("I120").Value = $ws.Range("I121").Value2
$ws.Range("J120").Value = 50
$ws.Range("K120").Value = 12000
$ws.Range("L120").Value = 13000
$ws.Range("M120").Value = 12600
$ws.Range("N120").Value = $ws.Range("N121").Value2
$ws.Range("O120").Value = $ws.Range("O121").Value2
$ws.Range("P120").Value = 210
$ws.Range("Q120").Value = $ws.Range("Q121").Value2
$ws.Range("R120").Value = $ws.Range("R121").Value2
